# David Silva section correction again
#
# "Lessons Learned" slide: the bullet about the dude who reads binary had a
# typo ("read in binary" -> "reads binary"). Fix the wording in place,
# leaving everything else (the rest of the bullet list, the other shapes,
# etc.) untouched.

$p = $ppt.ActivePresentation

# Find the slide that contains the bullet we need to fix instead of a
# hard-coded index, so the script is resilient to small structural changes.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -like "*read in binary*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$needle = "read in binary"
$startIdx = $fullText.IndexOf($needle)

# Replace just the faulty tail of the run ("read in binary" -> "reads
# binary"); PowerPoint keeps the untouched "The dude how " prefix as its own
# run and the retyped tail becomes a new run.
$target = $tr.Characters($startIdx + 1, $needle.Length)
$target.Text = "reads binary"
